$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.840.08"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.824.15"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.320"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "2.088.77"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "1.819.52"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.666"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "34.734.52"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "0.0₃0786"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("E29").Value = "  -5.96%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0549"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.698"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "91.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").Value = "1.337.93"
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  -3.08%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0523"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("D48").Value = "2.008.05"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0670"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
